$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "missing items" list gains one new product ("رباط ضغط 6سم"),
# which sorts in right after "جل صبار للبشره" (row 12) and before
# "فازلين هير تونك 200مل" (row 13). That pushes rows 13-15 down by one
# and creates a brand new row (item #10) that re-uses what used to be
# row 15's product. The totals row and the footer row each shift down
# one row as well, and the grand total grows by the new item's price.
# ------------------------------------------------------------------

# 1. Insert a fresh blank row right before the current totals row (16),
#    which shifts: old row16 (totals) -> 17, old row17 (footer) -> 18.
$ws.Rows("16:16").Insert()

# 2. Copy the formatting of row 15 (a normal item row) onto new row 16.
foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")) {
    $ws.Range($col + "15").Copy()
    $ws.Range($col + "16").PasteSpecial(-4122)  # xlPasteFormats
}
$ws.Application.CutCopyMode = $false

# 3. Match the row height used by the new row in the source workbook.
$ws.Rows("16:16").RowHeight = 25.5

# 4. Re-create the merged cells for the new row.
$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()

# 5. Row 16 (new item #10) gets the data that used to sit in row 15.
$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "مناديل مبلله كبيره"
$ws.Range("H16").Value = "5:0"
$ws.Range("L16").Value = "0"
$ws.Range("N16").Value = "30.00"
$ws.Range("P16").Value = "30.0000"
$ws.Range("Q16").Value = "1:0"

# 6. Row 15 (item #9) gets the data that used to sit in row 14.
$ws.Range("C15").Value = "كريم فريند لافلي الكبير"
$ws.Range("H15").Value = "20:0"
$ws.Range("N15").Value = "35.00"
$ws.Range("P15").Value = "35.0000"

# 7. Row 14 (item #8) gets the data that used to sit in row 13.
$ws.Range("C14").Value = "فازلين هير تونك 200مل"
$ws.Range("H14").Value = "0:0"
$ws.Range("N14").Value = "110.00"
$ws.Range("P14").Value = "110.0000"

# 8. Row 13 (item #7) becomes the brand new product.
$ws.Range("C13").Value = "رباط ضغط 6سم"
$ws.Range("H13").Value = "0:0"
$ws.Range("N13").Value = "15.00"
$ws.Range("P13").Value = "15.0000"

# 9. The grand total (now on row 17) grows by the new item's price.
$ws.Range("P17").Value = 314.5

# 10. The generated-on timestamp in the footer (now row 18) advances
#     two minutes, to match the re-uploaded version.
$ws.Range("A18").Value = "Saturday, 9 August, 2025 9:50 AM"
